$d = $word.ActiveDocument

# 1. Title line just under the document title: make it bold AND change its text.
$rTitleLine = $d.Content
$rTitleLine.Find.Execute("Document de référence de l'équipe de conception")
$rTitleLine.Bold = 1
$rTitleLine.Text = "Document des responsabilités principales de l'équipe de conception"

# 2. Remaining plain text replacements (old -> new), each unique within the document.
$pairs = @(
    @{old=" : ce document décrit les principales responsabilités de tous les membres de l'équipe de conception au Graphic Design Institute."; new=" : ce document décrit les principales responsabilités de tous les membres de l'équipe de conception du Graphic Design Institute."},

    @{old=" : collaborez avec d'autres concepteurs, développeurs et parties prenantes pour créer des conceptions de haute qualité qui répondent aux exigences du projet."; new=" : travailler en collaboration avec d'autres concepteurs, développeurs et parties prenantes pour créer des conceptions de haute qualité qui répondent aux exigences du projet,"},

    @{old="En outre, les concepteurs d'animation confirmés doivent travailler en collaboration avec d'autres concepteurs, développeurs et parties prenantes pour créer des conceptions de haute qualité qui répondent aux exigences du projet."; new="En outre, les concepteurs d'animation confirmés doivent travailler en collaboration avec d'autres concepteurs, développeurs et parties prenantes pour créer des conceptions de haute qualité qui répondent aux exigences du projet,"},

    @{old=" : créez des conceptions visuellement attrayantes qui sont conviviales, accessibles et réactives."; new=" : créer des conceptions visuellement attrayantes à la fois conviviales, accessibles et réactives,"},

    @{old="En outre, les concepteurs d'animation confirmés doivent créer des conceptions visuellement attrayantes à la fois conviviales, accessibles et réactives."; new="En outre, les concepteurs d'animation confirmés doivent créer des conceptions visuellement attrayantes à la fois conviviales, accessibles et réactives,"},

    @{old=" : communiquez efficacement avec les membres de l'équipe, les parties prenantes et les clients pour vous assurer que les exigences du projet sont remplies."; new=" : communiquer efficacement avec les membres de l'équipe, les parties prenantes et les clients pour s'assurer que les exigences du projet sont respectées,"},

    @{old="En outre, les concepteurs d'animation confirmés doivent communiquer efficacement avec les membres de l'équipe, les parties prenantes et les clients pour s'assurer que les exigences du projet sont bien respectées."; new="En outre, les concepteurs d'animation confirmés doivent communiquer efficacement avec les membres de l'équipe, les parties prenantes et les clients pour s'assurer que les exigences du projet sont bien respectées,"},

    @{old=" : Effectuez des recherches pour identifier les besoins, les préférences et les comportements des utilisateurs pour informer les décisions de conception."; new=" : effectuer des recherches pour identifier les besoins, les préférences et les comportements des utilisateurs afin d'informer les décisions de conception,"},

    @{old="En outre, les concepteurs d'animation confirmés doivent mener des recherches pour identifier les besoins, les préférences et les comportements des utilisateurs afin d'éclairer les décisions en matière de conception."; new="En outre, les concepteurs d'animation confirmés doivent mener des recherches pour identifier les besoins, les préférences et les comportements des utilisateurs afin d'éclairer les décisions en matière de conception,"},

    @{old="Test : effectuez des tests"; new="Test"},

    @{old=" d'utilisation pour vous assurer que les conceptions répondent aux besoins des utilisateurs et sont accessibles à tous les utilisateurs."; new=" : procéder à des tests de convivialité pour s'assurer que les conceptions répondent aux besoins des utilisateurs et sont accessibles à tous les utilisateurs,"},

    @{old="En outre, les concepteurs d'animation confirmés doivent procéder à des tests de convivialité pour s'assurer que les conceptions répondent aux besoins des utilisateurs et sont accessibles à tous les utilisateurs."; new="En outre, les concepteurs d'animation confirmés doivent procéder à des tests de convivialité pour s'assurer que les conceptions répondent aux besoins des utilisateurs et sont accessibles à tous les utilisateurs,"},

    @{old=" : Créez et gérez la documentation de conception, notamment les spécifications de conception, les repères de style et les modèles de conception."; new=" : créer et tenir à jour la documentation de conception, y compris les spécifications de conception, les guides de style et les modèles de conception,"},

    @{old="Développement"; new="Développement professionnel"},

    @{old=" professionnel : restez à jour avec les dernières tendances, outils et technologies de conception pour améliorer la qualité et l'efficacité de la conception."; new=" : se tenir informé des dernières tendances, outils et technologies de conception pour améliorer la qualité et l'efficacité de la conception,"},

    @{old="Leadership :"; new="Leadership"},

    @{old=" dirigez l'équipe de conception et fournissez des conseils aux concepteurs juniors."; new=" : diriger l'équipe de conception et fournir des conseils aux concepteurs débutants,"}
)

foreach ($p in $pairs) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($p.old, $true, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)
    Write-Output ($ok.ToString() + " :: " + $p.old)
}
